# Added support for custom organization:
# Insert 5 new rows into the "IDENTIFIEUR" list (column AL) right after the
# existing "BT16" entry (row 45) and before the start of the BT-code block
# (old row 46), then populate them with the new placeholder/organization
# entries: four "-" separators followed by a new "BT16" truck entry. All
# rows that used to start at row 46 shift down by 5 (to row 51+) keeping
# their original content and order intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 46:83 down to 51:88, inserting 5 blank rows at 46.
$ws.Range("A46:A50").EntireRow.Insert()

# Populate the newly inserted rows.
$ws.Range("AL46").Value = "-"
$ws.Range("AL47").Value = "-"
$ws.Range("AL48").Value = "-"
$ws.Range("AL49").Value = "-"
$ws.Range("AL50").Value = "BT16"

# Restore the cursor/selection roughly where the author left it when saving.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 38
$null = $ws.Range("AN38").Select()
